$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.048.30"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").Value = "1.650.38"
$ws.Range("E3").Value = "  -0.58%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.36"
$ws.Range("E5").Value = "  +0.28%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5260"
$ws.Range("E6").Value = "  +2.06%  "
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2598"
$ws.Range("E8").Value = "  -1.70%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06315"
$ws.Range("E9").Value = "  +0.72%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.36"
$ws.Range("E10").Value = "  -2.09%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07803"
$ws.Range("E11").Value = "  +0.66%  "
$ws.Range("E12").Value = "  +0.63%  "
$ws.Range("D13").Value = "1.647.75"
$ws.Range("E13").Value = "  -1.01%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5486"
$ws.Range("E14").Value = "  +0.68%  "
$ws.Range("D15").Value = "0.0₅8197"
$ws.Range("E15").Value = "  +1.18%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "65.36"
$ws.Range("E16").Value = "  +0.72%  "
$ws.Range("D17").Value = "26.068.29"
$ws.Range("E17").Value = "  -0.55%  "
$ws.Range("E18").Value = "  -0.16%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "4.573"
$ws.Range("E19").Value = "  -0.92%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "191.05"
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("E21").Value = "  -0.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.040"
$ws.Range("E22").Value = "  +0.50%  "
$ws.Range("E23").Value = "  -0.17%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "142.49"
$ws.Range("E24").Value = "  +2.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1236"
$ws.Range("E25").Value = "  +1.05%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.224"
$ws.Range("E26").Value = "  -0.47%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.04"
$ws.Range("E27").Value = "  -0.75%  "
$ws.Range("E28").Value = "  -0.53%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.05809"
$ws.Range("E29").Value = "  -2.58%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.272"
$ws.Range("E30").Value = "  +0.05%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.538"
$ws.Range("E31").Value = "  -0.75%  "
$ws.Range("E32").Value = "  -0.01%  "
$ws.Range("E33").Value = "  -0.46%  "
$ws.Range("E34").Value = "  -0.57%  "
$ws.Range("B35").Value = "MXToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.777"
$ws.Range("E35").Value = "  +0.36%  "
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9439"
$ws.Range("E36").Value = "  -2.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5730"
$ws.Range("E37").Value = "  +0.82%  "
$ws.Range("E38").Value = "  +0.85%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.8443"
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.747"
$ws.Range("E40").Value = "  -5.10%  "
$ws.Range("E41").Value = "  -0.07%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "103.48"
$ws.Range("E42").Value = "  +3.21%  "
$ws.Range("D43").Value = "1.027.32"
$ws.Range("E43").Value = "  +1.60%  "
$ws.Range("D44").Value = "1.794.96"
$ws.Range("E44").Value = "  -0.37%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "57.01"
$ws.Range("E45").Value = "  +0.53%  "
$ws.Range("E46").Value = "  -0.78%  "
$ws.Range("E47").Value = "  +3.01%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.850"
$ws.Range("E48").Value = "  -2.28%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05145"
$ws.Range("E49").Value = "  -0.46%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.466"
$ws.Range("E50").Value = "  +1.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.09635"
$ws.Range("E51").Value = "  -0.56%  "
